# Badge Order Tracking.xlsx - "Add files via upload"
#
# The uploaded workbook adds one more row to the PCB order-tracking table:
# a "Duty on PCB Boards" line item (row 25), which was previously a blank
# spacer row between the table and the rows below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the new order row (row 25) with the duty/customs line item.
$ws.Range("A25").Value = "Duty on PCB Boards"
$ws.Range("B25").Value = "Duty on Main PCB Board Order"
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = "DHL"
$ws.Range("E25").Value = 313.13
$ws.Range("F25").Value = "kevin"
$ws.Range("H25").Value = 44820
$ws.Range("I25").Value = 44820

# Arrival/received date columns use the same centered "d-mmm" date format
# as the other rows in the table (e.g. H9, H17, H19).
$ws.Range("H25").NumberFormat = "d-mmm"
$ws.Range("I25").NumberFormat = "d-mmm"

# Matches the author's final on-screen selection after the edit.
$ws.Range("A26").Select() | Out-Null
